$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($row, $col, $value) {
    $cell = $ws.Cells.Item($row, $col)
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

function Set-PlainValue($row, $col, $value) {
    $ws.Cells.Item($row, $col).Value = $value
}

# Row 2
Set-TextValue 2 4 '60.693.53'
Set-PlainValue 2 5 '  -2.84%  '
# Row 3
Set-TextValue 3 4 '2.901.53'
Set-PlainValue 3 5 '  -3.89%  '
# Row 4
Set-TextValue 4 4 '1.00'
Set-PlainValue 4 5 '  +0.00%  '
# Row 5
Set-TextValue 5 4 '585.70'
Set-PlainValue 5 5 '  -1.81%  '
# Row 6
Set-TextValue 6 4 '146.60'
Set-PlainValue 6 5 '  -2.49%  '
# Row 7
Set-PlainValue 7 5 '  +0.09%  '
# Row 8
Set-PlainValue 8 5 '  -3.08%  '
# Row 9
Set-TextValue 9 4 '2.901.37'
Set-PlainValue 9 5 '  -3.80%  '
# Row 10
Set-TextValue 10 4 '6.68'
Set-PlainValue 10 5 '  +4.59%  '
# Row 11
Set-TextValue 11 4 '0.143'
Set-PlainValue 11 5 '  -4.37%  '
# Row 13
Set-PlainValue 13 5 '  -3.71%  '
# Row 14
Set-TextValue 14 4 '33.82'
Set-PlainValue 14 5 '  -2.24%  '
# Row 16
Set-TextValue 16 4 '3.383.04'
Set-PlainValue 16 5 '  -3.88%  '
# Row 17
Set-TextValue 17 4 '6.81'
Set-PlainValue 17 5 '  -2.76%  '
# Row 18
Set-TextValue 18 4 '60.633.50'
Set-PlainValue 18 5 '  -2.95%  '
# Row 19
Set-TextValue 19 4 '2.899.56'
Set-PlainValue 19 5 '  -4.04%  '
# Row 20
Set-TextValue 20 4 '425.44'
Set-PlainValue 20 5 '  -5.28%  '
# Row 21
Set-PlainValue 21 5 '  -3.97%  '
# Row 22
Set-PlainValue 22 5 '  -2.74%  '
# Row 23
Set-PlainValue 23 5 '  -4.77%  '
# Row 24
Set-TextValue 24 4 '80.32'
Set-PlainValue 24 5 '  -2.52%  '
# Row 25
Set-TextValue 25 4 '11.05'
Set-PlainValue 25 5 '  +1.26%  '
# Row 26
Set-TextValue 26 4 '2.24'
Set-PlainValue 26 5 '  -0.30%  '
# Row 27
Set-PlainValue 27 5 '  -1.77%  '
# Row 28
Set-PlainValue 28 5 '  +0.06%  '
# Row 29
Set-TextValue 29 4 '7.31'
Set-PlainValue 29 5 '  +0.18%  '
# Row 30
Set-TextValue 30 4 '2.21'
Set-PlainValue 30 5 '  +3.06%  '
# Row 31
Set-PlainValue 31 5 '  -0.19%  '
# Row 32
Set-TextValue 32 4 '2.62'
Set-PlainValue 32 5 '  -3.62%  '
# Row 33
Set-TextValue 33 4 '26.47'
Set-PlainValue 33 5 '  -3.84%  '
# Row 34
Set-PlainValue 34 5 '  -2.53%  '
# Row 35
Set-TextValue 35 4 '0.0₃0835'
Set-PlainValue 35 5 '  -1.95%  '
# Row 36
Set-PlainValue 36 5 '  -2.12%  '
# Row 37
Set-TextValue 37 4 '5.65'
Set-PlainValue 37 5 '  -3.39%  '
# Row 38
Set-PlainValue 38 2 'Stacks'
Set-PlainValue 38 3 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
Set-TextValue 38 4 '2.02'
Set-PlainValue 38 5 '  -2.20%  '
# Row 39
Set-PlainValue 39 2 'OKB'
Set-PlainValue 39 3 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
Set-TextValue 39 4 '49.19'
Set-PlainValue 39 5 '  -1.79%  '
# Row 40
Set-TextValue 40 4 '2.95'
Set-PlainValue 40 5 '  -3.35%  '
# Row 41
Set-TextValue 41 4 '0.124'
Set-PlainValue 41 5 '  -0.41%  '
# Row 42
Set-TextValue 42 4 '8.72'
Set-PlainValue 42 5 '  -3.60%  '
# Row 43
Set-TextValue 43 4 '0.294'
Set-PlainValue 43 5 '  +2.84%  '
# Row 44
Set-TextValue 44 4 '41.87'
Set-PlainValue 44 5 '  +4.20%  '
# Row 45
Set-PlainValue 45 5 '  -2.04%  '
# Row 46
Set-TextValue 46 4 '371.60'
Set-PlainValue 46 5 '  -4.87%  '
# Row 47
Set-PlainValue 47 2 'Maker'
Set-PlainValue 47 3 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
Set-TextValue 47 4 '2.653.47'
Set-PlainValue 47 5 '  -3.03%  '
# Row 48
Set-PlainValue 48 2 'Monero'
Set-PlainValue 48 3 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
Set-TextValue 48 4 '133.18'
Set-PlainValue 48 5 '  +0.36%  '
# Row 49
Set-TextValue 49 4 '25.46'
Set-PlainValue 49 5 '  +7.15%  '
# Row 51
Set-PlainValue 51 5 '  -1.22%  '
